# Update the "Förändrad" (Changed) date column (C) for all data rows
# from 45174 (2023-09-05) to 45175 (2023-09-06).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45174) {
        $cell.Value = 45175
    }
}
